$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "[Paroma%Deb%xref no email%1,       Kaniz\u2010e\u2010%Zannat%xref no email%1,       Shiny%Talukder%xref no email%1,       Amirul Huda%Bhuiyan%xref no email%1,       Md. Shariful Alam%Jilani%xref no email%1,       K. M.%Saif\u2010Ur\u2010Rahman%xref no email%1]"
$ws.Range("E3").Value = "[Hussein N.%Ali%NULL%1,                          Sherko S.%Niranji%sherko.subhan@garmian.edu.krd%1,                          Sirwan M. A.%Al‐Jaf%NULL%2,                          Sirwan M. A.%Al‐Jaf%NULL%0]"
$ws.Range("E4").Value = "[Ana Valesca Fernandes Gilson%Silva%NULL%1,                          Diego%Menezes%NULL%1,                          Filipe Romero Rebello%Moreira%NULL%1,                          Octávio Alcântara%Torres%NULL%1,                          Paula Luize Camargos%Fonseca%NULL%1,                          Rennan Garcias%Moreira%NULL%1,                          Hugo José%Alves%NULL%1,                          Vivian Ribeiro%Alves%NULL%1,                          Tânia Maria de Resende%Amaral%NULL%1,                          Adriano Neves%Coelho%NULL%1,                          Júlia Maria%Saraiva Duarte%NULL%1,                          Augusto Viana%da Rocha%NULL%1,                          Luiz Gonzaga Paula%de Almeida%NULL%1,                          João Locke Ferreira%de Araújo%NULL%1,                          Hilton Soares%de Oliveira%NULL%1,                          Nova Jersey Cláudio%de Oliveira%NULL%1,                          Camila%Zolini%NULL%1,                          Jôsy Hubner%de Sousa%NULL%1,                          Elizângela Gonçalves%de Souza%NULL%1,                          Rafael Marques%de Souza%NULL%1,                          Luciana de Lima%Ferreira%NULL%1,                          Alexandra%Lehmkuhl Gerber%NULL%1,                          Ana Paula de Campos%Guimarães%NULL%1,                          Paulo Henrique Silva%Maia%NULL%1,                          Fernanda Martins%Marim%NULL%1,                          Lucyene%Miguita%NULL%1,                          Cristiane Campos%Monteiro%NULL%1,                          Tuffi Saliba%Neto%NULL%1,                          Fabrícia Soares Freire%Pugêdo%NULL%1,                          Daniel Costa%Queiroz%NULL%1,                          Damares Nigia Alborguetti Cuzzuol%Queiroz%NULL%1,                          Luciana Cunha%Resende-Moreira%NULL%1,                          Franciele Martins%Santos%NULL%1,                          Erika Fernanda Carlos%Souza%NULL%1,                          Carolina Moreira%Voloch%NULL%1,                          Ana Tereza%Vasconcelos%NULL%1,                          Renato Santana%de Aguiar%NULL%1,                          Renan Pedra%de Souza%NULL%1]"
